# "Add files via upload" — 이정원 adds a new Work-Log entry documenting the
# first pass of Code Smell clean-up on Member.java, then leaves her own
# sheet as the active tab/selection when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("이정원")

# New row of data (row 5): 할일 / 내용 / 할당일 / 완성일 / 결과 / 문제점
$ws.Range("A5").Value = "Code Smell 제거"
$ws.Range("B5").Value = "Member.java 의 코드 스멜 1차 제거"
$ws.Range("C5").Value = "2019-05-17"
$ws.Range("D5").Value = "2019-05-17"
$ws.Range("E5").Value = "개선된 코드 Member.java를 git commit함"

# Match the row height used for the other filled-in rows on this sheet.
$ws.Range("A5:F5").EntireRow.RowHeight = 35

# Page setup was touched for this sheet (portrait, paper size 9 = A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# This sheet becomes the active tab, with E6 as the selected cell.
$ws.Activate()
$ws.Range("E6").Select() | Out-Null
